# Update of league bases (Lithuania A Lyga) - re-order of duplicate-date match
# rows. The rows of the same round/date got reshuffled: for each affected
# group of rows, everything except the "id" column (A) is rotated among the
# rows in the group, i.e. the match-level data (match id, teams, odds, P/L...)
# moves to a different row while the running row index in column A stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AB$r1")
    $rng2 = $ws.Range("B$r2`:AB$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Simple pairwise swaps (B:AB only - column A "id" stays fixed per row)
Swap-Rows 50 51
Swap-Rows 89 90
Swap-Rows 100 101
Swap-Rows 117 118
Swap-Rows 136 137

# Rows 102/103/104 rotate in a 3-cycle:
#   new(102) = old(103), new(103) = old(104), new(104) = old(102)
$rng102 = $ws.Range("B102:AB102")
$rng103 = $ws.Range("B103:AB103")
$rng104 = $ws.Range("B104:AB104")

$v102 = $rng102.Value2
$v103 = $rng103.Value2
$v104 = $rng104.Value2

$rng102.Value = $v103
$rng103.Value = $v104
$rng104.Value = $v102

Write-Output "Row data reshuffled successfully"
